$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.305692434310913
$ws.Range("B1").Value = 3.924838066101074
$ws.Range("C1").Value = 3.813581943511963
$ws.Range("D1").Value = 3.056463479995728
$ws.Range("E1").Value = 1.043466091156006
